$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($r, $id, $cat, $type, $qty, $amount, $date, $time) {
    $ws.Cells.Item($r, 1).Value = $id
    $ws.Cells.Item($r, 2).Value = $cat
    $ws.Cells.Item($r, 3).Value = $type
    $ws.Cells.Item($r, 4).Value = $qty
    $ws.Cells.Item($r, 5).Value = $amount
    # The literal date text (e.g. "2024-09-20") looks like a real date to
    # Excel's input parser and would otherwise be auto-converted into a
    # date serial number. Enter it with a trailing space (blocks the
    # date-parsing heuristic) then normalize it back to the exact text via
    # a TRIM() formula copied in as a value - this keeps the cell a plain
    # text cell with the default style (no NumberFormat/style changes
    # leak into styles.xml).
    $ws.Cells.Item($r, 6).Value = ($date + " ")
    $helper = $ws.Range("ZZ1")
    $helper.Formula = "=TRIM(" + $ws.Cells.Item($r, 6).Address($false, $false) + ")"
    $helper.Copy() | Out-Null
    $ws.Cells.Item($r, 6).PasteSpecial(-4163) | Out-Null
    $helper.Clear() | Out-Null

    # Time-like text (e.g. "20:29:33") is not auto-converted by the
    # engine, so it can be assigned directly.
    $ws.Cells.Item($r, 7).Value = $time
}

Set-Row 13 "22ed083d-1f01-4059-be79-16814d3b979e" "In" "One Face" 100 10 "2024-09-20" "20:29:33"
Set-Row 14 "e5b9c353-440b-47a1-b36e-ca67e8d62331" "In" "Duable Face" 100 15 "2024-09-20" "20:29:51"
Set-Row 15 "70979251-9e54-4aab-baf7-d2be5fb0d6ac" "Waste" "paper" 100 0 "2024-09-20" "20:43:29"
Set-Row 16 "4214537e-1c4a-4fa1-9e5f-27df97d6fbee" "In" "Duable Face" 100 15 "2024-09-20" "20:44:19"
